# ---------------------------------------------------------------------------
# Línea 141: horarios completos del 30/12/2025
#
# Updates the "horarios-141-completo.xlsx" workbook with the latest scraped
# arrival data for the three sheets: LP1912, LP1912-215 and 6203-6173.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ===========================================================================
# Sheet 1: LP1912
# ===========================================================================
$ws1 = $wb.Worksheets.Item("LP1912")

# --- new "Fecha" header column (G1), matching the bold/bordered header style
$ws1.Range("F1").Copy()
$ws1.Range("G1").PasteSpecial($xlPasteFormats)
$ws1.Range("G1").Value2 = "Fecha"

# --- refresh the "last update" / "total rows" meta rows
$ws1.Range("A2").Value2 = "Última actualización: 29/12/2025 22:53:01"
$ws1.Range("A3").Value2 = "Total filas: 8"

# --- create the (empty / unstyled) placeholder cells for G2:G4 and A5:A9,
#     reusing the plain, style-less formatting already used by A4
$ws1.Range("A4").Copy()
$ws1.Range("G2:G4").PasteSpecial($xlPasteFormats)
$ws1.Range("A5:A9").PasteSpecial($xlPasteFormats)

# --- new data rows scraped on 29/12/2025
$ws1.Range("B5").Value2 = "22:52:50"
$ws1.Range("C5").Value2 = "23:06"
$ws1.Range("D5").Value2 = "16_SANTA ANA"
$ws1.Range("E5").Value2 = 14
$ws1.Range("F5").Value2 = "LP1912"
$ws1.Range("G5").Value2 = "29/12/2025"

$ws1.Range("B6").Value2 = "22:52:50"
$ws1.Range("C6").Value2 = "23:07"
$ws1.Range("D6").Value2 = "14X44_ABASTO"
$ws1.Range("E6").Value2 = 15
$ws1.Range("F6").Value2 = "LP1912"
$ws1.Range("G6").Value2 = "29/12/2025"

$ws1.Range("B7").Value2 = "22:52:50"
$ws1.Range("C7").Value2 = "23:15"
$ws1.Range("D7").Value2 = "23_HERNANDEZ"
$ws1.Range("E7").Value2 = 23
$ws1.Range("F7").Value2 = "LP1912"
$ws1.Range("G7").Value2 = "29/12/2025"

$ws1.Range("B8").Value2 = "22:52:50"
$ws1.Range("C8").Value2 = "23:49"
$ws1.Range("D8").Value2 = "16_SANTA ANA"
$ws1.Range("E8").Value2 = 57
$ws1.Range("F8").Value2 = "LP1912"
$ws1.Range("G8").Value2 = "29/12/2025"

$ws1.Range("B9").Value2 = "22:52:50"
$ws1.Range("C9").Value2 = "23:51"
$ws1.Range("D9").Value2 = "215_ALUAR"
$ws1.Range("E9").Value2 = 59
$ws1.Range("F9").Value2 = "LP1912"
$ws1.Range("G9").Value2 = "29/12/2025"

# ===========================================================================
# Sheet 2: LP1912-215
# ===========================================================================
$ws2 = $wb.Worksheets.Item("LP1912-215")

# --- promote A1 to the bold header style and add the new header columns
$ws1.Range("A1:F1").Copy()
$ws2.Range("A1:F1").PasteSpecial($xlPasteFormats)

$ws1.Range("G1").Copy()
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$ws2.Range("B1").Value2 = "Fecha"
$ws2.Range("C1").Value2 = "Hora_Scrap"
$ws2.Range("D1").Value2 = "Hora_Llegada"
$ws2.Range("E1").Value2 = "Línea"
$ws2.Range("F1").Value2 = "Minutos"
$ws2.Range("G1").Value2 = "Parada"

# --- refresh the "last update" / "total rows" meta rows
$ws2.Range("A2").Value2 = "Última actualización: 29/12/2025 22:53:01"
$ws2.Range("A3").Value2 = "Total filas: 3"

# --- create the (empty / unstyled) placeholder cells for B2:G3 and A4,
#     reusing the plain, style-less formatting already used by A4 on LP1912
$ws1.Range("A4").Copy()
$ws2.Range("B2:G2").PasteSpecial($xlPasteFormats)
$ws2.Range("B3:G3").PasteSpecial($xlPasteFormats)
$ws2.Range("A4").PasteSpecial($xlPasteFormats)

# --- new data row scraped on 29/12/2025
$ws2.Range("B4").Value2 = "29/12/2025"
$ws2.Range("C4").Value2 = "22:52:50"
$ws2.Range("D4").Value2 = "23:51"
$ws2.Range("E4").Value2 = "215_ALUAR"
$ws2.Range("F4").Value2 = 59
$ws2.Range("G4").Value2 = "LP1912"

# ===========================================================================
# Sheet 3: 6203-6173
# ===========================================================================
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- promote A1 to the bold header style and add the new header columns
$ws1.Range("A1:F1").Copy()
$ws3.Range("A1:F1").PasteSpecial($xlPasteFormats)

$ws1.Range("G1").Copy()
$ws3.Range("G1").PasteSpecial($xlPasteFormats)

$ws3.Range("B1").Value2 = "Fecha"
$ws3.Range("C1").Value2 = "Hora_Scrap"
$ws3.Range("D1").Value2 = "Hora_Llegada"
$ws3.Range("E1").Value2 = "Línea"
$ws3.Range("F1").Value2 = "Minutos"
$ws3.Range("G1").Value2 = "Parada"

# --- refresh the "last update" / "total rows" meta rows
$ws3.Range("A2").Value2 = "Última actualización: 29/12/2025 22:53:01"
$ws3.Range("A3").Value2 = "Total filas: 3"

# --- create the (empty / unstyled) placeholder cells for B2:G3 and A4,
#     reusing the plain, style-less formatting already used by A4 on LP1912
$ws1.Range("A4").Copy()
$ws3.Range("B2:G2").PasteSpecial($xlPasteFormats)
$ws3.Range("B3:G3").PasteSpecial($xlPasteFormats)
$ws3.Range("A4").PasteSpecial($xlPasteFormats)

# --- new data row scraped on 29/12/2025
$ws3.Range("B4").Value2 = "29/12/2025"
$ws3.Range("C4").Value2 = "22:53:01"
$ws3.Range("D4").Value2 = "23:08"
$ws3.Range("E4").Value2 = "215A_LA PLATA"
$ws3.Range("F4").Value2 = 15
$ws3.Range("G4").Value2 = "L6173"
